$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Metadata")
$ws2 = $wb.Worksheets.Item("Include from RxNorm")

# --- Sheet1 "Metadata" updates ---

# Version 0.1.6 -> 0.1.7
$ws1.Range("B3").Value = "0.1.7"

# Status active -> draft
$ws1.Range("B6").Value = "draft"

# Date updated
$ws1.Range("B8").Value = "2024-08-23T10:17:11-05:00"

# Contact text updated (row 10)
$ws1.Range("B10").Value = "The Medical College of Wisconsin, Inc. and the National Marrow Donor Program (http://www.cibmtr.org)"

# Second contact (row 11) updated in place (used to be a duplicate "No display for ContactDetail")
$ws1.Range("B11").Value = "Bob Milius (bmilius@nmdp.org)"

# Insert a new row 12 for "Jurisdiction" (empty value), pushing Description/Purpose/Copyright/Immutable down by one
$ws1.Rows.Item(12).Insert()
$ws1.Range("A14:B14").Copy()
$ws1.Range("A12:B12").PasteSpecial(-4122)
$ws1.Range("A12").Value = "Jurisdiction"
